$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A-H (rows 2-9) to be stored as text, matching the source
# data (account / contract numbers must not be reinterpreted as numbers).
$ws.Range("A2:H9").NumberFormat = "@"

# Data rows 2-8 (new content), columns A-K
$data = @(
    @("NOUBAIL MOHAMMED", "IR801997", "007400000313200019604463", "KHOURIBGA ZELLAKA", "AWB", "Direction régionale", "035/TES/AV1", "mensuelle", 1000, 100, 900),
    @("NOUBAIL MOHAMMED", "IR801997", "007400000313200019604463", "KHOURIBGA ZELLAKA", "AWB", "Direction régionale", "035/TES/AV1", "mensuelle", 4000, 400, 3600),
    @("NASIRI HASNAA", "", "546576878798989898090090", "", "CIH", "Logement de fonction", "905/LF/TADLA OUARDIGHA ZAYANE", "mensuelle", 9999.99, 999.99, 9000),
    @("ZERNAKH ABDELLAH", "IB19558", "145101211406073828000084", "MARRAKECH BENI MELLAL", "BP", "Point de vente", "052/FKIH BEN SALEH/AV1", "mensuelle", 12000, 0, 12000),
    @("MOHAMED BADRANE", "I83603", "225400000805987601012173", "KHOURIBGA", "CA", "Point de vente", "605/KHOURIBGA NAHDA", "mensuelle", 7500, 375, 7125),
    @("NOUBAIL MOUNTASSIR", "Q251990", "007400000313200019604463", "KHOURIBGA ZELLAKA", "AWB", "Direction régionale", "905/TADLA OUARDIGHA ZAYANE", "mensuelle", 6750, 675, 6075),
    @("NOUBAIL MOHAMMED", "IR801997", "007400000313200019604463", "KHOURIBGA ZELLAKA", "AWB", "Direction régionale", "905/TADLA OUARDIGHA ZAYANE", "mensuelle", 6750, 675, 6075)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $row++
}

# Totals row, now at row 9
$ws.Cells.Item(9, 1).Value = " "
$ws.Cells.Item(9, 2).Value = " "
$ws.Cells.Item(9, 3).Value = " "
$ws.Cells.Item(9, 4).Value = " "
$ws.Cells.Item(9, 5).Value = " "
$ws.Cells.Item(9, 6).Value = " "
$ws.Cells.Item(9, 7).Value = " "
$ws.Cells.Item(9, 8).Value = " "
$ws.Cells.Item(9, 9).Value = 47999.99
$ws.Cells.Item(9, 10).Value = 3224.99
$ws.Cells.Item(9, 11).Value = 44775
